$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.591.41'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.718.04'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.26%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9983'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.46'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9986'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4927'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.80%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2601'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.98%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06203'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.725.58'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07000'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6056'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.477'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.75'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9986'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.439.29'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9986'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007142'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.946.31'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.408'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.512'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.50%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.074'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '137.79'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.88%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.403'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.740'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '105.57'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.919'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07947'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.637'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.39%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04509'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.11%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9973'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6239'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9343'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.999'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.408'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9982'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01510'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.22'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.512'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3835'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.907'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05373'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.764'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.16'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '51.45'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.223'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.43%  '
